$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed trace-report pull returned fewer events (7 instead of 9) and a
# later completion timestamp. Two of the previously-reported "Placed Actual"
# car events (rows 6 & 7: HRTX541048 / CRDX15008) are gone from the new pull,
# and the old "Arrive In-Transit / KANSAS CITY / HLINKC" event for HRTX541043
# has been superseded by a later "Departure / HUTCHINSON / HKCKDE" event.

# Remove the two rows that no longer appear in the refreshed pull.
# (This shifts rows 8-13 up to 6-11 automatically, carrying their
# formatting/styles with them.)
$ws.Rows("6:7").Delete() | Out-Null

# The HRTX541043 row (now row 9 after the shift above) needs its event
# details updated to the latest status from the refreshed pull.
$ws.Range("C9").Value = "HUTCHINSON"
$ws.Range("F9").Value = 22
$ws.Range("G9").Value = 845
$ws.Range("H9").Value = "Departure"
$ws.Range("I9").Value = "HKCKDE"

# Update the summary header cells to match the refreshed pull.
$ws.Range("A1").Value = "Description unknown, completed 06/22/2023 11:07:25 EDT, by WPJTOWN1.The search returned: 7 events."
$ws.Range("A2").Value = "1 On Hand"

# Keep the selection / filter database in sync with the now-smaller table.
$ws.Range("K5:K11").Select() | Out-Null

$names = $wb.Names
$fdb = $names.Item("_xlnm._FilterDatabase")
$fdb.RefersTo = '=Test_format_trace!$A$4:$O$11'
